$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 2
    4  = 2
    5  = 2
    6  = 0
    7  = 2
    8  = 1
    9  = 2
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 2
    15 = 3
    16 = 2
    17 = 2
    18 = 0
    19 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
